$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The H2:H13 "is_active" cells currently hold a =TRUE() boolean formula
# (cached numeric value 1). They need to become the literal text string
# "TRUE" instead. A direct .Value/.Formula assignment of "TRUE" gets
# auto-detected as a Boolean by Excel, so instead we type it with a
# leading apostrophe (forcing literal text) into a scratch cell, copy
# it, and paste-special (values only) onto the target cell - this keeps
# the destination's existing style/number-format intact while changing
# its stored type to a shared text string. The scratch cell is fully
# cleared afterwards so it leaves no trace.
for ($r = 2; $r -le 13; $r++) {
    $scratch = $ws.Cells.Item($r, 7)
    $scratch.Value = "'TRUE"
    $scratch.Copy()
    $target = $ws.Cells.Item($r, 8)
    $target.PasteSpecial(-4163)
    $scratch.Clear()
}

$excel.CutCopyMode = $false

# Update the active selection to match the edited range
$ws.Range("H2:H13").Select()
